$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1229.5
$ws.Range("I33").Value = 1247.2858
$ws.Range("K33").Value = 1247.2858
$ws.Range("M33").Value = -1018.2858

$ws.Range("H62").Value = 5698
$ws.Range("I62").Value = 6149.5
$ws.Range("K62").Value = 6149.5
$ws.Range("M62").Value = -5525.5

$ws.Range("H65").Value = 5698
$ws.Range("I65").Value = 6149.5
$ws.Range("K65").Value = 30747.5
$ws.Range("M65").Value = -27627.5

$ws.Range("H96").Value = 241.85715
$ws.Range("J96").Value = 300
$ws.Range("L96").Value = 900
$ws.Range("N96").Value = -3646

$ws.Range("H97").Value = 485.7143
$ws.Range("J97").Value = 466.66666
$ws.Range("L97").Value = 1399.99998
$ws.Range("N97").Value = -2391.99998

$ws.Range("H100").Value = 5183.1577
$ws.Range("I100").Value = 2354.75
$ws.Range("J100").Value = 7240.1816
$ws.Range("K100").Value = 2354.75
$ws.Range("L100").Value = 7240.1816
$ws.Range("M100").Value = -1813.75
$ws.Range("N100").Value = -8322.1816

$ws.Range("H112").Value = 1607.9454
$ws.Range("J112").Value = 1610.4073
$ws.Range("L112").Value = 4831.2219
$ws.Range("N112").Value = -7047.2219

$ws.Range("H137").Value = 1894.5555
$ws.Range("I137").Value = 1285.7894
$ws.Range("J137").Value = 3340.375
$ws.Range("K137").Value = 3857.3682
$ws.Range("L137").Value = 10021.125
$ws.Range("M137").Value = -1307.3682
$ws.Range("N137").Value = -15121.125

$ws.Range("H138").Value = 2336.6216
$ws.Range("I138").Value = 1158.76
$ws.Range("J138").Value = 4790.5
$ws.Range("K138").Value = 3476.28
$ws.Range("L138").Value = 14371.5
$ws.Range("M138").Value = 1663.72
$ws.Range("N138").Value = -24651.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10058.242
$ws.Range("I32").Value = 7055.5625
$ws.Range("K32").Value = 7055.5625
$ws.Range("M32").Value = -6768.5625

$ws.Range("H45").Value = 6306.8
$ws.Range("I45").Value = 4630.8
$ws.Range("J45").Value = 7982.8
$ws.Range("K45").Value = 4630.8
$ws.Range("L45").Value = 7982.8
$ws.Range("M45").Value = -4253.8
$ws.Range("N45").Value = -8736.799999999999

$ws.Range("H53").Value = 8000
$ws.Range("I53").Value = 8000
$ws.Range("K53").Value = 8000
$ws.Range("M53").Value = -7318

$ws.Range("H61").Value = 3762.6553
$ws.Range("J61").Value = 5226.385
$ws.Range("L61").Value = 5226.385
$ws.Range("N61").Value = -5650.385

$ws.Range("H63").Value = 3826.3462
$ws.Range("I63").Value = 3738.652
$ws.Range("J63").Value = 4498.6665
$ws.Range("K63").Value = 3738.652
$ws.Range("L63").Value = 4498.6665
$ws.Range("M63").Value = -3052.652
$ws.Range("N63").Value = -5870.6665

$ws.Range("H66").Value = 3826.3462
$ws.Range("I66").Value = 3738.652
$ws.Range("J66").Value = 4498.6665
$ws.Range("K66").Value = 18693.26
$ws.Range("L66").Value = 22493.3325
$ws.Range("M66").Value = -15261.26
$ws.Range("N66").Value = -29357.3325

$ws.Range("H102").Value = 4151.6665
$ws.Range("I102").Value = 4151.6665
$ws.Range("K102").Value = 4151.6665
$ws.Range("M102").Value = -2529.6665

$ws.Range("H133").Value = 68717.414

$ws.Range("H136").Value = 3762.6553
$ws.Range("J136").Value = 5226.385
$ws.Range("L136").Value = 15679.155
$ws.Range("N136").Value = -20779.155

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2290.7144
$ws.Range("I105").Value = 1394.0625
$ws.Range("K105").Value = 1394.0625
$ws.Range("M105").Value = 352.9375

$ws.Range("H134").Value = 4579.25
$ws.Range("I134").Value = 3005.8667
$ws.Range("K134").Value = 9017.6001
$ws.Range("M134").Value = -6482.6001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2934.4119
$ws.Range("I58").Value = 1878.7778
$ws.Range("K58").Value = 1878.7778
$ws.Range("M58").Value = -1675.7778

$ws.Range("H105").Value = 55562730
$ws.Range("J105").Value = 125018870
$ws.Range("L105").Value = 125018870
$ws.Range("N105").Value = -125022364

$ws.Range("H135").Value = 69908.63

$ws.Range("H136").Value = 2934.4119
$ws.Range("I136").Value = 1878.7778
$ws.Range("K136").Value = 5636.3334
$ws.Range("M136").Value = -3086.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 619.7273
$ws.Range("I8").Value = 619.7273
$ws.Range("K8").Value = 1859.1819
$ws.Range("M8").Value = -1720.1819

$ws.Range("H23").Value = 1116.8462
$ws.Range("J23").Value = 416
$ws.Range("L23").Value = 1248
$ws.Range("N23").Value = -1718

$ws.Range("H38").Value = 170.90909
$ws.Range("J38").Value = 307.5
$ws.Range("L38").Value = 922.5
$ws.Range("N38").Value = -1616.5

$ws.Range("H80").Value = 4625.273
$ws.Range("I80").Value = 4102.5
$ws.Range("J80").Value = 4924
$ws.Range("K80").Value = 12307.5
$ws.Range("L80").Value = 14772
$ws.Range("M80").Value = -11371.5
$ws.Range("N80").Value = -16644

$ws.Range("H83").Value = 4625.273
$ws.Range("I83").Value = 4102.5
$ws.Range("J83").Value = 4924
$ws.Range("K83").Value = 36922.5
$ws.Range("L83").Value = 44316
$ws.Range("M83").Value = -32242.5
$ws.Range("N83").Value = -53676

$ws.Range("H107").Value = 15151772
$ws.Range("I107").Value = 184.5
$ws.Range("J107").Value = 16666931
$ws.Range("K107").Value = 553.5
$ws.Range("L107").Value = 50000793
$ws.Range("M107").Value = 1366.5
$ws.Range("N107").Value = -50004633

$ws.Range("H133").Value = 1666.6666
$ws.Range("I133").Value = 1000
$ws.Range("K133").Value = 3000
$ws.Range("M133").Value = 2060

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3976.5938
$ws.Range("I132").Value = 3137.3684
$ws.Range("K132").Value = 9412.1052
$ws.Range("M132").Value = -6882.1052

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 10368.477
$ws.Range("I40").Value = 13068
$ws.Range("K40").Value = 13068
$ws.Range("M40").Value = -12932

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1994
$ws.Range("I96").Value = 1994
$ws.Range("J96").Value = 1994
$ws.Range("K96").Value = 1994
$ws.Range("L96").Value = 1994
$ws.Range("M96").Value = -621
$ws.Range("N96").Value = -4740

$ws.Range("H100").Value = 649.3570999999999
$ws.Range("I100").Value = 649.3570999999999
$ws.Range("K100").Value = 1298.7142
$ws.Range("M100").Value = -757.7141999999999
